$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Q21 (row 22) benchmark count after re-benching on updated nightly.
$ws.Range("C22").Value = 587

# Formulas in D22, H22, I22, D24, H24 recalc automatically.

# Reflect the cell selection left by the author when saving.
$ws.Range("C23").Select()
